# Add a "password" column for the Create-account feature and a new player row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# A new "password" column is inserted at D, pushing the existing "coin"
# column from D to E and "char_class" from E to F.
$coinHeader = $ws.Range("D1").Value()
$charClassHeader = $ws.Range("E1").Value()

$ws.Range("F1").Value = $charClassHeader
$ws.Range("E1").Value = $coinHeader
$ws.Range("D1").Value = "password"

# --- New data row (row 3) for player "Loc Le" ---
$ws.Range("A3").Value = "Loc Le"
$ws.Range("B3").Value = "ltloc05lumia520@gmail.com"
$ws.Range("C3").Value = "Loc Le.png"
$ws.Range("D3").Value = "khongbiethehe"

# --- Column width for the newly added password column (D) ---
$ws.Columns.Item(4).ColumnWidth = 12.3

# --- Selection moves to G7, matching the authored workbook state ---
$ws.Range("G7").Select() | Out-Null
